$d = $word.ActiveDocument

# Curly single quotes used throughout the document.
$lq = [char]0x2018
$rq = [char]0x2019

# 'image' <BinData>  ->  'imageUrl' <string>
$find    = $lq + "image" + $rq + " <BinData>"
$replace = $lq + "imageUrl" + $rq + " <string>"
[void]$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)

# Word leaves its "_GoBack" bookmark marking the site of the last edit -
# i.e. right after the newly typed "string" (before the trailing '>').
$full = $d.Content.Text
$titIdx = $full.IndexOf("titlul")
$bookmarkPos = $titIdx + 46
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
